# Remove the "Email" feature: delete column C (Email) entirely, shifting
# the remaining columns (Reason, Amount, Account Number, Account Name,
# Bank Name) one to the left, and update the sample row with new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column C ("Email") - this shifts D:H left to C:G and updates
# the used-range dimension automatically (A1:H2 -> A1:G2).
$ws.Columns("C").Delete()

# Update the data row (row 2) with the new values from the commit.
$ws.Range("A2").Value = "2025-05-18 21:57"
$ws.Range("D2").Value = 5000
$ws.Range("E2").Value = 9055301016
